$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$found = $d.Content
$ok = $found.Find.Execute("Docente(s) Responsável(eis)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($ok) {
    # Find the paragraph that contains the matched range.
    $targetIndex = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $found.Start -and $p.Range.End -ge $found.End) {
            $targetIndex = $i
            break
        }
    }

    $targetPara = $d.Paragraphs.Item($targetIndex)

    # Insert a new empty paragraph right after it.
    $targetPara.Range.InsertParagraphAfter()

    # The newly created paragraph now sits right after the target one;
    # give it the bulleted-list style and the new text.
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "5701460 - Antonio Iacono"
    $newPara.Range.Style = "ListBullet"
}
